$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds date-serial values. Every populated data row
# (rows 2-43) currently has the value 45787 (2025-05-10) and should be
# bumped to 45788 (2025-05-11).
for ($r = 2; $r -le 43; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45787) {
        $cell.Value2 = 45788
    }
}
